$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "26.547.12", "  +0.71%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.729.70", "  +0.64%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.9995", "  +0.00%  ")
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "245.28", "  +2.78%  ")
    ,@(6, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.9998", "  -0.05%  ")
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.4805", "  +1.57%  ")
    ,@(8, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2673", "  +1.35%  ")
    ,@(9, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.06224", "  +0.21%  ")
    ,@(10, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.727.94", "  +0.61%  ")
    ,@(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07160", "  +1.46%  ")
    ,@(12, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "15.71", "  +2.50%  ")
    ,@(13, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6163", "  +4.09%  ")
    ,@(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.543", "  +2.94%  ")
    ,@(15, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "77.18", "  +1.26%  ")
    ,@(16, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.9998", "  -0.06%  ")
    ,@(17, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "26.548.13", "  +0.77%  ")
    ,@(18, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "0.9997", "  -0.07%  ")
    ,@(19, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000006977", "  +2.47%  ")
    ,@(20, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "11.66", "  +0.93%  ")
    ,@(21, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "1.950.29", "  +0.68%  ")
    ,@(22, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "4.533", "  -0.33%  ")
    ,@(23, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.933", "  +1.89%  ")
    ,@(24, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "5.286", "  -0.65%  ")
    ,@(25, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "136.56", "  +1.18%  ")
    ,@(26, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "15.35", "  +0.71%  ")
    ,@(27, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.800", "  +2.62%  ")
    ,@(28, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.405", "  +0.03%  ")
    ,@(29, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "106.78", "  -1.20%  ")
    ,@(30, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "3.991", "  -0.38%  ")
    ,@(31, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.07992", "  +3.23%  ")
    ,@(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "3.722", "  +1.02%  ")
    ,@(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.04594", "  +3.68%  ")
    ,@(34, "Frax", "https://coinranking.com/coin/KfWtaeV1W+frax-frax", "0.9995", "  -0.05%  ")
    ,@(35, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.617", "  +0.16%  ")
    ,@(36, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.6363", "  +2.72%  ")
    ,@(37, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "0.9964", "  +1.84%  ")
    ,@(38, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.9206", "  -1.37%  ")
    ,@(39, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "2.091", "  +9.16%  ")
    ,@(40, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.404", "  -0.32%  ")
    ,@(41, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "105.00", "  -7.72%  ")
    ,@(42, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.003", "  +0.30%  ")
    ,@(43, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01507", "  +1.97%  ")
    ,@(44, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "5.585", "  +4.89%  ")
    ,@(45, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.3898", "  +2.04%  ")
    ,@(46, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "6.997", "  +11.20%  ")
    ,@(47, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1185", "  +1.29%  ")
    ,@(48, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05340", "  +1.07%  ")
    ,@(49, "Elrond", "https://coinranking.com/coin/omwkOTglq+elrond-egld", "30.98", "  +1.92%  ")
    ,@(50, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "7.886", "  +2.35%  ")
    ,@(51, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "1.265", "  +4.00%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
